$d = $word.ActiveDocument
$newText = "Ημερομηνίες παρατήρησης για τον αστερισμό του Cygnus: 10-19 Αυγούστου, 9-18 Σεπτεμβρίου, 8-17 Οκτωβρίου"

# Locate every paragraph whose text begins with the old "2018 ..." observation-date
# sentence (there are 4 of them in the document) and rebuild it as a single plain run.
$targets = New-Object System.Collections.ArrayList
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.StartsWith("2018 ") -and $t -like "*30 Οκτωβρίου-8 Νοεμβρίου*") {
        [void]$targets.Add($i)
    }
}

# Process from last to first so earlier paragraph indices stay valid.
for ($k = $targets.Count - 1; $k -ge 0; $k--) {
    $idx = $targets[$k]
    $p = $d.Paragraphs($idx)
    $r = $p.Range
    $r.End = $r.End - 1
    $r.Delete()

    $r2 = $p.Range
    $r2.End = $r2.End - 1
    $r2.Text = $newText
}
